$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.212.92"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.82%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.853.45"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +1.50%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.33%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "313.61"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.30%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4639"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3713"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07279"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.88%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.8869"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.08%  "

$ws.Range("E11").Value = "  +1.75%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07856"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.52%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.917.34"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +5.15%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.388"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.93%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "6.517"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.46%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "90.88"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("E17").Value = "  -0.34%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000008914"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("E19").Value = "  -0.30%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "14.70"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "27.238.05"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.84%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.083"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.27%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "10.51"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.150.51"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +4.65%  "

$ws.Range("E25").Value = "  +5.71%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "151.43"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -1.03%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "18.40"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.045"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.04%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "115.86"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.054"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.49%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.08809"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.76%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.141"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +6.28%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.7675"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +5.34%  "

$ws.Range("E34").Value = "  +3.08%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.508"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +1.50%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.731"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +10.18%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.111"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +4.03%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01940"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.38%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.05217"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("E40").Value = "  -0.40%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "7.027"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.95%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.5116"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.05%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.1627"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.18%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "8.440"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +3.38%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.4794"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.02%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "10.35"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +1.65%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "102.73"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.643"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.83%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.06202"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "65.58"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.31%  "
